$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "da"
$ws.Range("C5").Value = "dsa"

$ws.Range("C5").Select()
